$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Helper: write a value as plain text without leaving behind a
# permanent number-format/style change on the cell (Excel would
# otherwise auto-detect strings that look like numbers, e.g. "311.91",
# and store them as numeric cells).
function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Update price (D) and volume-change (E) columns for each affected coin row.
$ws.Range("D2").Value = "27.553.92"
$ws.Range("E2").Value = "  -0.37%  "

$ws.Range("D3").Value = "1.835.84"
$ws.Range("E3").Value = "  -0.51%  "

$ws.Range("E4").Value = "  -0.08%  "

Set-TextValue $ws.Range("D5") "311.91"
$ws.Range("E5").Value = "  -0.22%  "

Set-TextValue $ws.Range("D7") "0.4277"
$ws.Range("E7").Value = "  -0.01%  "

Set-TextValue $ws.Range("D8") "0.3653"
$ws.Range("E8").Value = "  +0.58%  "

Set-TextValue $ws.Range("D9") "0.07269"
$ws.Range("E9").Value = "  -0.81%  "

Set-TextValue $ws.Range("D10") "0.8647"
$ws.Range("E10").Value = "  -1.36%  "

Set-TextValue $ws.Range("D11") "20.67"
$ws.Range("E11").Value = "  +0.23%  "

$ws.Range("D12").Value = "1.826.42"
$ws.Range("E12").Value = "  -2.50%  "

Set-TextValue $ws.Range("D13") "5.438"
$ws.Range("E13").Value = "  +1.57%  "

Set-TextValue $ws.Range("D14") "6.522"
$ws.Range("E14").Value = "  +0.05%  "

Set-TextValue $ws.Range("D15") "0.06955"
$ws.Range("E15").Value = "  +0.05%  "

$ws.Range("E16").Value = "  -0.09%  "

Set-TextValue $ws.Range("D17") "80.62"
$ws.Range("E17").Value = "  +1.32%  "

Set-TextValue $ws.Range("D18") "0.000008914"
$ws.Range("E18").Value = "  -0.28%  "

$ws.Range("E19").Value = "  -0.19%  "

Set-TextValue $ws.Range("D20") "15.43"
$ws.Range("E20").Value = "  +0.58%  "

$ws.Range("D21").Value = "27.560.22"
$ws.Range("E21").Value = "  -0.91%  "

Set-TextValue $ws.Range("D22") "5.158"
$ws.Range("E22").Value = "  +3.46%  "

Set-TextValue $ws.Range("D23") "10.85"
$ws.Range("E23").Value = "  +5.21%  "

$ws.Range("D24").Value = "2.068.07"
$ws.Range("E24").Value = "  -2.95%  "

Set-TextValue $ws.Range("D25") "1.995"
$ws.Range("E25").Value = "  +0.33%  "

Set-TextValue $ws.Range("D26") "154.76"
$ws.Range("E26").Value = "  -0.47%  "

Set-TextValue $ws.Range("D27") "18.92"
$ws.Range("E27").Value = "  +2.25%  "

Set-TextValue $ws.Range("D28") "5.150"
$ws.Range("E28").Value = "  -1.09%  "

Set-TextValue $ws.Range("D29") "114.20"
$ws.Range("E29").Value = "  -4.52%  "

Set-TextValue $ws.Range("D30") "1.826"
$ws.Range("E30").Value = "  -2.31%  "

Set-TextValue $ws.Range("D31") "0.08868"
$ws.Range("E31").Value = "  -0.06%  "

Set-TextValue $ws.Range("D32") "2.999"
$ws.Range("E32").Value = "  +1.18%  "

Set-TextValue $ws.Range("D33") "0.7489"
$ws.Range("E33").Value = "  -0.71%  "

Set-TextValue $ws.Range("D34") "4.543"
$ws.Range("E34").Value = "  +0.56%  "

Set-TextValue $ws.Range("D35") "1.132"
$ws.Range("E35").Value = "  +0.18%  "

$ws.Range("E36").Value = "  -0.06%  "

Set-TextValue $ws.Range("D37") "1.096"
$ws.Range("E37").Value = "  -1.01%  "

Set-TextValue $ws.Range("D38") "0.05324"
$ws.Range("E38").Value = "  -2.09%  "

Set-TextValue $ws.Range("D39") "0.01937"
$ws.Range("E39").Value = "  +0.28%  "

$ws.Range("E40").Value = "  -0.70%  "

Set-TextValue $ws.Range("D41") "0.5083"
$ws.Range("E41").Value = "  +0.27%  "

$ws.Range("E42").Value = "  -0.65%  "

Set-TextValue $ws.Range("D43") "6.488"
$ws.Range("E43").Value = "  -1.55%  "

Set-TextValue $ws.Range("D44") "8.334"
$ws.Range("E44").Value = "  -0.62%  "

Set-TextValue $ws.Range("D45") "10.46"
$ws.Range("E45").Value = "  +0.98%  "

Set-TextValue $ws.Range("D48") "0.4673"
$ws.Range("E48").Value = "  +0.52%  "

Set-TextValue $ws.Range("D49") "0.9997"
$ws.Range("E49").Value = "  -0.14%  "

Set-TextValue $ws.Range("D50") "1.620"
$ws.Range("E50").Value = "  -1.07%  "

Set-TextValue $ws.Range("D51") "63.65"
$ws.Range("E51").Value = "  -1.69%  "

# Row 46/47: Cronos and Quant swapped order, with updated prices and changes.
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Range("D46") "105.46"
$ws.Range("E46").Value = "  -0.61%  "

$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D47") "0.06472"
$ws.Range("E47").Value = "  -1.15%  "

